$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 470, pushing old rows 470-555 down to 472-557.
$ws.Rows("470:471").Insert()

# New row 470 data
$r = 470
$ws.Cells.Item($r,1).Value2  = 10
$ws.Cells.Item($r,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item($r,3).Value2  = "La Araucanía"
$ws.Cells.Item($r,4).Value2  = 45180
$ws.Cells.Item($r,5).Value2  = 9
$ws.Cells.Item($r,6).Value2  = 100112017
$ws.Cells.Item($r,7).Value2  = "Apio"
$ws.Cells.Item($r,8).Value2  = "Americana (o)"
$ws.Cells.Item($r,9).Value2  = "Primera"
$ws.Cells.Item($r,10).Value2 = 300
$ws.Cells.Item($r,11).Value2 = 8000
$ws.Cells.Item($r,12).Value2 = 8000
$ws.Cells.Item($r,13).Value2 = 8000
$ws.Cells.Item($r,14).Value2 = "`$/caja 8 unidades"
$ws.Cells.Item($r,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item($r,16).Value2 = 8000
$ws.Cells.Item($r,17).Value2 = 1
$ws.Cells.Item($r,18).Value2 = "Hortaliza"

# New row 471 data
$r = 471
$ws.Cells.Item($r,1).Value2  = 10
$ws.Cells.Item($r,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item($r,3).Value2  = "La Araucanía"
$ws.Cells.Item($r,4).Value2  = 45180
$ws.Cells.Item($r,5).Value2  = 9
$ws.Cells.Item($r,6).Value2  = 100112017
$ws.Cells.Item($r,7).Value2  = "Apio"
$ws.Cells.Item($r,8).Value2  = "Americana (o)"
$ws.Cells.Item($r,9).Value2  = "Primera"
$ws.Cells.Item($r,10).Value2 = 210
$ws.Cells.Item($r,11).Value2 = 7000
$ws.Cells.Item($r,12).Value2 = 8000
$ws.Cells.Item($r,13).Value2 = 7714
$ws.Cells.Item($r,14).Value2 = "`$/docena de matas"
$ws.Cells.Item($r,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item($r,16).Value2 = 1286
$ws.Cells.Item($r,17).Value2 = 6
$ws.Cells.Item($r,18).Value2 = "Hortaliza"
